$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.584.43"
$ws.Cells.Item(2, 5).Value = "  -1.54%  "
$ws.Cells.Item(3, 4).Value = "1.666.19"
$ws.Cells.Item(3, 5).Value = "  -3.21%  "
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).Value = "214.90"
$ws.Cells.Item(5, 5).Value = "  -1.64%  "
$ws.Cells.Item(6, 4).Value = "0.509"
$ws.Cells.Item(6, 5).Value = "  -2.36%  "
$ws.Cells.Item(7, 5).Value = "  -0.06%  "
$ws.Cells.Item(8, 4).Value = "23.80"
$ws.Cells.Item(8, 5).Value = "  -1.61%  "
$ws.Cells.Item(9, 4).Value = "0.261"
$ws.Cells.Item(9, 5).Value = "  -0.76%  "
$ws.Cells.Item(10, 5).Value = "  -1.68%  "
$ws.Cells.Item(11, 4).Value = "0.0879"
$ws.Cells.Item(11, 5).Value = "  -2.25%  "
$ws.Cells.Item(12, 4).Value = "1.902.87"
$ws.Cells.Item(12, 5).Value = "  -3.13%  "
$ws.Cells.Item(13, 4).Value = "1.651.26"
$ws.Cells.Item(13, 5).Value = "  -4.06%  "
$ws.Cells.Item(14, 4).Value = "4.13"
$ws.Cells.Item(14, 5).Value = "  -3.35%  "
$ws.Cells.Item(15, 4).Value = "0.559"
$ws.Cells.Item(15, 5).Value = "  -0.27%  "
$ws.Cells.Item(16, 4).Value = "66.22"
$ws.Cells.Item(16, 5).Value = "  -1.85%  "
$ws.Cells.Item(17, 4).Value = "27.579.33"
$ws.Cells.Item(17, 5).Value = "  -1.38%  "
$ws.Cells.Item(18, 4).Value = "242.68"
$ws.Cells.Item(18, 5).Value = "  +0.29%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0728"
$ws.Cells.Item(19, 5).Value = "  -3.31%  "
$ws.Cells.Item(20, 4).Value = "7.63"
$ws.Cells.Item(20, 5).Value = "  -4.10%  "
$ws.Cells.Item(21, 5).Value = "  -0.13%  "
$ws.Cells.Item(22, 5).Value = "  -3.05%  "
$ws.Cells.Item(23, 4).Value = "9.28"
$ws.Cells.Item(23, 5).Value = "  -3.80%  "
$ws.Cells.Item(24, 5).Value = "  -4.70%  "
$ws.Cells.Item(25, 4).Value = "146.98"
$ws.Cells.Item(25, 5).Value = "  -1.11%  "
$ws.Cells.Item(26, 4).Value = "7.19"
$ws.Cells.Item(26, 5).Value = "  -3.98%  "
$ws.Cells.Item(27, 4).Value = "16.43"
$ws.Cells.Item(27, 5).Value = "  -1.55%  "
$ws.Cells.Item(28, 5).Value = "  -0.14%  "
$ws.Cells.Item(29, 4).Value = "0.111"
$ws.Cells.Item(29, 5).Value = "  -2.35%  "
$ws.Cells.Item(30, 5).Value = "  +2.94%  "
$ws.Cells.Item(31, 5).Value = "  -1.39%  "
$ws.Cells.Item(32, 4).Value = "3.34"
$ws.Cells.Item(32, 5).Value = "  -2.49%  "
$ws.Cells.Item(33, 4).Value = "1.465.83"
$ws.Cells.Item(33, 5).Value = "  -1.44%  "
$ws.Cells.Item(34, 4).Value = "3.11"
$ws.Cells.Item(34, 5).Value = "  -4.89%  "
$ws.Cells.Item(35, 4).Value = "1.56"
$ws.Cells.Item(35, 5).Value = "  -5.31%  "
$ws.Cells.Item(36, 5).Value = "  -1.42%  "
$ws.Cells.Item(37, 4).Value = "0.926"
$ws.Cells.Item(37, 5).Value = "  -2.81%  "
$ws.Cells.Item(38, 5).Value = "  -1.36%  "
$ws.Cells.Item(39, 4).Value = "0.573"
$ws.Cells.Item(39, 5).Value = "  -5.52%  "
$ws.Cells.Item(40, 4).Value = "69.42"
$ws.Cells.Item(40, 5).Value = "  -1.66%  "
$ws.Cells.Item(41, 5).Value = "  -5.16%  "
$ws.Cells.Item(42, 5).Value = "  -0.11%  "
$ws.Cells.Item(43, 2).Value = "MXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(43, 4).Value = "2.22"
$ws.Cells.Item(43, 5).Value = "  -3.09%  "
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).Value = "5.40"
$ws.Cells.Item(44, 5).Value = "  -7.28%  "
$ws.Cells.Item(45, 4).Value = "1.809.92"
$ws.Cells.Item(45, 5).Value = "  -3.03%  "
$ws.Cells.Item(46, 5).Value = "  -1.28%  "
$ws.Cells.Item(47, 5).Value = "  -2.57%  "
$ws.Cells.Item(48, 4).Value = "89.27"
$ws.Cells.Item(48, 5).Value = "  -1.69%  "
$ws.Cells.Item(49, 5).Value = "  -3.12%  "
$ws.Cells.Item(50, 5).Value = "  -1.85%  "
$ws.Cells.Item(51, 4).Value = "7.85"
$ws.Cells.Item(51, 5).Value = "  -4.77%  "
